$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 needs to become the text value "1" (stored as a shared string, not a
# number) while keeping its existing style (s="23"). A plain
# `$ws.Range("B11").Value = "1"` gets auto-coerced to the number 1 by Excel's
# usual "smart" type detection, and forcing NumberFormat="@" first (or
# prefixing with an apostrophe) changes the cell's style (adds a text /
# quote-prefix format), which we don't want either.
#
# The trick: stage the literal text "1" in a scratch cell that is explicitly
# formatted as Text, then Copy / PasteSpecial *values only* into B11. Paste
# Special (values) carries over the cell's data type (text) without touching
# the destination's existing formatting/style.

$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

# Clean up the scratch cell so it doesn't leave a stray value behind.
$scratch.Clear()
